$wb = $excel.ActiveWorkbook

$oldGuid = "36a770d4-3d99-47de-96b9-0c595a2532dc"
$newGuid = "3984c643-9cb2-409a-9e87-cf9af1c1a39e"
$oldHash = "a6f301ab3a890d2eced44f7466a74c725d42beb3"
$newHash = "1afb124ed452b6cdafa2c5d9b29e114bef0afc74"

# Hyperlink "Address" targets are unchanged (still point at the original commit URL) -
# only the visible display text is updated to the new report's filename.
$hyperlinkUrl = "https://github.com/OpenLocalizationTestOrg/oltest/blob/88bb4afd36dc36f86e2f01fb0331623ae842dfdf/e2e/$oldGuid.md"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkUrl, [System.Type]::Missing, [System.Type]::Missing, "e2e\$newGuid.md")
$wsOverview.Range("G2").Value = "2016-08-12 19:11:50"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $hyperlinkUrl, [System.Type]::Missing, [System.Type]::Missing, "$newGuid.md")
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-12 19:11:43"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $hyperlinkUrl, [System.Type]::Missing, [System.Type]::Missing, "$newGuid.md")
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
# de-de!H2 shares the same underlying string as Overview!G2 ("Latest HO Xliff
# Generate Date" / "Latest Handoff DateTime" both held "2016-08-12 19:11:21"),
# so it also becomes "2016-08-12 19:11:50".
$wsDeDe.Range("H2").Value = "2016-08-12 19:11:50"
